$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Europe columns (D:G) entirely; Japan/Saudi/US columns shift left.
$ws.Range("D1:G1").EntireColumn.Delete()

# Update the re-worded question labels in column A.
$ws.Range("A5").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""
$ws.Range("A6").Value = "Would support a global movement to tackle CC, tax millionaires,`n and fund LICs (either petition, demonstrate, strike, or donate)"

# Refresh the data values (rows 2-9, columns B:K) with the final figures.
$ws.Range("B2").Value = 0.556007222541788
$ws.Range("C2").Value = 0.392706735271895
$ws.Range("D2").Value = 0.356105750421195
$ws.Range("E2").Value = 0.508432438341536
$ws.Range("F2").Value = 0.451519556659133
$ws.Range("G2").Value = 0.677827112481047
$ws.Range("H2").Value = 0.743498584988873
$ws.Range("I2").Value = 0.503741905167505
$ws.Range("J2").Value = 0.637096655800029
$ws.Range("K2").Value = 0.373993319775201
$ws.Range("B3").Value = 0.50312291438834
$ws.Range("C3").Value = 0.316198486427608
$ws.Range("D3").Value = 0.34215953895184
$ws.Range("E3").Value = 0.455068700122256
$ws.Range("F3").Value = 0.306824263084082
$ws.Range("G3").Value = 0.671355992475847
$ws.Range("H3").Value = 0.73716679715381
$ws.Range("I3").Value = 0.455679906453223
$ws.Range("J3").Value = 0.566496834134014
$ws.Range("K3").Value = 0.339067925103543
$ws.Range("B4").Value = 0.680881448179833
$ws.Range("C4").Value = 0.616567982061628
$ws.Range("D4").Value = 0.743644347389163
$ws.Range("E4").Value = 0.814701212857562
$ws.Range("F4").Value = 0.757048871605567
$ws.Range("G4").Value = 0.713280127381035
$ws.Range("H4").Value = 0.703520370125625
$ws.Range("I4").Value = 0.671270631778761
$ws.Range("J4").Value = 0.776836935461012
$ws.Range("K4").Value = 0.425661149175785
$ws.Range("B5").Value = 0.609601586795904
$ws.Range("C5").Value = 0.459499864440356
$ws.Range("D5").Value = 0.529969608967456
$ws.Range("E5").Value = 0.619684579205792
$ws.Range("F5").Value = 0.581343297592584
$ws.Range("G5").Value = 0.756844881931732
$ws.Range("H5").Value = 0.822839088961434
$ws.Range("I5").Value = 0.529112697724995
$ws.Range("J5").Value = 0.580830159607851
$ws.Range("K5").Value = 0.389019482696932
$ws.Range("B6").Value = 0.675595447215337
$ws.Range("C6").Value = 0.523930159271177
$ws.Range("D6").Value = 0.433349195600366
$ws.Range("E6").Value = 0.696851480613757
$ws.Range("F6").Value = 0.583790255087382
$ws.Range("G6").Value = 0.727098526374066
$ws.Range("H6").Value = 0.741985444624183
$ws.Range("I6").Value = 0.641824096726743
$ws.Range("J6").Value = 0.834461320073758
$ws.Range("K6").Value = 0.474126518973143
$ws.Range("B7").Value = 0.364717906507653
$ws.Range("C7").Value = 0.301242387158432
$ws.Range("D7").Value = 0.163327499246366
$ws.Range("E7").Value = 0.313179598308858
$ws.Range("F7").Value = 0.218483573122562
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 0.283095421815601
$ws.Range("J7").Value = 0.52122515690493
$ws.Range("K7").Value = 0.27084855688435
$ws.Range("B8").Value = 0.347853243460036
$ws.Range("C8").Value = 0.262527011404327
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 0.311745273790548
$ws.Range("J8").Value = 0.450524011973634
$ws.Range("K8").Value = 0.173904958357855
$ws.Range("B9").Value = 0.410626908494325
$ws.Range("C9").Value = 0.393121093071121
$ws.Range("D9").Value = 0.282669471326983
$ws.Range("E9").Value = 0.376571407830385
$ws.Range("F9").Value = 0.32600471502799
$ws.Range("G9").Value = 0.662536527281334
$ws.Range("H9").Value = 0.698489352155402
$ws.Range("I9").Value = 0.362475436951784
$ws.Range("J9").Value = 0.535384805366787
$ws.Range("K9").Value = 0.265617828927838
